$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "아이언디바이스" listing (previously on row 12, between "넥스트바이오메디컬"
# and "유라클") moved up the demand-forecast calendar: a new row is inserted
# right after "에이치이엠파마" (row 2) with an updated demand-forecast date,
# and its old row lower in the sheet is removed so every other listing keeps
# its original relative order.

# 1) Insert a fresh blank row at row 3 - everything from row 3 down shifts to row 4 down.
$ws.Rows.Item(3).Insert()

# 2) The old "아이언디바이스" row, originally row 12, is now row 13 after the insert.
#    Remove it so the listing appears only once, in its new spot.
$ws.Rows.Item(13).Delete()

# 3) Populate the newly inserted row 3 with "아이언디바이스"'s data, carrying over
#    the same price range / confirmed price / offering amount / underwriter,
#    but with the refreshed demand-forecast date.
$ws.Range("A3").Value = "아이언디바이스"
$ws.Range("B3").Value = "2024.08.19~08.23"
$ws.Range("C3").Value = "4,900~5,700"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 14700
$ws.Range("F3").Value = "대신증권"
